# RPA datasets push 2024-04-13
# Insert a new IPO row ("에스오에스랩") into the "02_38커뮤니케이션(최근일자기준)"
# sheet right after "노브랜드" (row 3), pushing all following rows down by
# one, and drop what is now the trailing duplicate row so the table keeps
# its original extent (21 rows incl. header).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Shift rows 4:21 down to 5:22 by inserting a blank row at row 4.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row with the new listing's data.
$ws.Cells.Item(4, 1).Value = "에스오에스랩"
$ws.Cells.Item(4, 2).Value = "2024.04.30~05.08"
$ws.Cells.Item(4, 3).Value = "7,500~9,000"
$ws.Cells.Item(4, 4).Value = "-"
$ws.Cells.Item(4, 5).Value = 15000
$ws.Cells.Item(4, 6).Value = "한국투자증권"

# The insert pushed the old last data row (row 21, "오상헬스케어") down to
# row 22 -- remove it so the sheet returns to its original 21-row extent.
$ws.Rows("22:22").Delete()
